$p = $ppt.ActivePresentation
$f = $p.Fonts
Write-Output "Fonts.Count: $($f.Count)"
for ($i=1; $i -le $f.Count; $i++) {
  $ff = $f.Item($i)
  Write-Output "$i : $($ff.Name)"
}
